# Resultados_SmartScore — add xiaoyi's second submission (row 30), duplicated
# from row 29, and fix row 29's SmartScore cells to be real numbers instead
# of text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# SmartScore columns: numeric on row 29 after this edit; row 30 keeps the
# original text-formatted values (as row 29 had before the edit).
$scoreCols = @("I", "L", "O", "R", "U", "X", "AA", "AD", "AG")
$scoreVals = @(0.591, 0.532, 0.518, 0.607, 0.52, 0.443, 0.718, 0.705, 0.673)
# Text renditions exactly as they were originally stored (row 29 kept a
# trailing zero on "0.520" that [string] formatting of the double would drop).
$scoreText = @("0.591", "0.532", "0.518", "0.607", "0.520", "0.443", "0.718", "0.705", "0.673")

# Every other column gets copied verbatim from row 29 into the new row 30.
$copyCols = @("C", "D", "E", "G", "H", "J", "K", "M", "N", "P", "Q", "S", "T", "V", "W", "Y", "Z", "AB", "AC", "AE", "AF", "AH")

# --- Row 30: duplicate of row 29 -------------------------------------------
foreach ($col in $copyCols) {
    $ws.Range($col + "30").Value2 = $ws.Range($col + "29").Value2
}

# B29 is an (empty-string) text cell — a bare quote keeps B30 text/empty too
# instead of leaving the cell completely blank/absent.
$ws.Range("B30").Value2 = "'"

# SmartScore text values (quote-prefixed so they stay text, matching the
# original formatting of these cells before row 29 was converted to numbers).
for ($i = 0; $i -lt $scoreCols.Count; $i++) {
    $col = $scoreCols[$i]
    $ws.Range($col + "30").Value2 = "'" + $scoreText[$i]
}

# New participant id / submission timestamp for the duplicated row.
$ws.Range("A30").Value2 = "xiaoyi_20251202_134617"
$ws.Range("F30").Value2 = "2025-12-02 13:46:18"

# --- Row 29: SmartScore text values become real numbers --------------------
for ($i = 0; $i -lt $scoreCols.Count; $i++) {
    $col = $scoreCols[$i]
    $ws.Range($col + "29").Value2 = $scoreVals[$i]
}

# Keep the new row's height the same as every other (default) row — writing
# the multi-line "Pesos" JSON text would otherwise auto-expand it.
$ws.Rows("30:30").RowHeight = $ws.Rows("29:29").RowHeight

Write-Host "SmartScore row duplicated and row 29 scores converted to numbers."
